$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in the sheet
$lastRow = $ws.UsedRange.Rows.Count

# Column C holds the "Förändrad" (Changed) date values.
# All data rows (2..lastRow) currently store serial date 46061 (2026-02-08)
# and must be bumped by one day to 46062 (2026-02-09).
$range = $ws.Range("C2:C$lastRow")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -ne $null) {
        $cell.Value2 = $cell.Value2 + 1
    }
}
